$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.964.98'
$ws.Range("D3").Value = '1.878.47'
$ws.Range("E3").Value = '  -2.00%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9974'
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.73'
$ws.Range("E5").Value = '  -4.89%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9976'
$ws.Range("E6").Value = '  -0.18%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5000'
$ws.Range("E7").Value = '  -2.98%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '44.50'
$ws.Range("E8").Value = '  -3.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2926'
$ws.Range("E9").Value = '  -1.73%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06616'
$ws.Range("E10").Value = '  -3.49%  '
$ws.Range("D11").Value = '1.877.88'
$ws.Range("E11").Value = '  -2.01%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '16.74'
$ws.Range("E12").Value = '  -4.49%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07255'
$ws.Range("E13").Value = '  -1.11%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6675'
$ws.Range("E14").Value = '  -3.47%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '86.17'
$ws.Range("E15").Value = '  -1.94%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.862'
$ws.Range("E16").Value = '  -1.20%  '
$ws.Range("D17").Value = '29.940.49'
$ws.Range("E17").Value = '  -1.12%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007896'
$ws.Range("E18").Value = '  -3.12%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9977'
$ws.Range("E19").Value = '  -0.21%  '
$ws.Range("E20").Value = '  -2.40%  '
$ws.Range("D21").Value = '2.119.77'
$ws.Range("E21").Value = '  -2.04%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9974'
$ws.Range("E22").Value = '  -0.06%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.752'
$ws.Range("E23").Value = '  -2.32%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.643'
$ws.Range("E24").Value = '  -2.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.057'
$ws.Range("E25").Value = '  -1.50%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '149.09'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '139.34'
$ws.Range("E27").Value = '  -0.09%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.96'
$ws.Range("E28").Value = '  -1.88%  '
$ws.Range("E29").Value = '  -5.11%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.384'
$ws.Range("E30").Value = '  +0.32%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.179'
$ws.Range("E31").Value = '  -2.52%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08793'
$ws.Range("E32").Value = '  -0.76%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.944'
$ws.Range("E33").Value = '  -2.20%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05051'
$ws.Range("E34").Value = '  -1.75%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7154'
$ws.Range("E35").Value = '  -0.77%  '
$ws.Range("E36").Value = '  -4.48%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.663'
$ws.Range("E37").Value = '  -0.78%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.692'
$ws.Range("E38").Value = '  -4.72%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01749'
$ws.Range("E39").Value = '  +2.65%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.183'
$ws.Range("E40").Value = '  -5.43%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9297'
$ws.Range("E41").Value = '  -4.96%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4259'
$ws.Range("E42").Value = '  -1.58%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.791'
$ws.Range("E43").Value = '  -5.50%  '
$ws.Range("E44").Value = '  -0.20%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '101.71'
$ws.Range("E45").Value = '  -4.29%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.446'
$ws.Range("E46").Value = '  -3.57%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1256'
$ws.Range("E47").Value = '  -1.78%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05654'
$ws.Range("E48").Value = '  -1.48%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '32.36'
$ws.Range("E49").Value = '  -3.25%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.3755'
$ws.Range("E50").Value = '  -2.35%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.218'
$ws.Range("E51").Value = '  -3.71%  '
